{"js": "// Replace each three-digit x one-digit multiplication expression with its\n// corresponding updated expression, one find/replace pair per cell.\nconst replacements = [\n  [\"581\u00d73=1743\", \"530\u00d73=1590\"],\n  [\"551\u00d73=1653\", \"639\u00d79=5751\"],\n  [\"411\u00d72=822\", \"514\u00d74=2056\"],\n  [\"789\u00d77=5523\", \"124\u00d73=372\"],\n  [\"960\u00d75=4800\", \"469\u00d79=4221\"],\n  [\"872\u00d79=7848\", \"806\u00d79=7254\"],\n  [\"263\u00d77=1841\", \"978\u00d77=6846\"],\n  [\"716\u00d77=5012\", \"278\u00d76=1668\"],\n  [\"369\u00d78=2952\", \"343\u00d74=1372\"],\n  [\"980\u00d78=7840\", \"336\u00d72=672\"],\n  [\"389\u00d78=3112\", \"245\u00d74=980\"],\n  [\"346\u00d74=1384\", \"392\u00d72=784\"],\n  [\"538\u00d74=2152\", \"825\u00d72=1650\"],\n  [\"179\u00d75=895\", \"974\u00d72=1948\"],\n  [\"554\u00d73=1662\", \"604\u00d72=1208\"],\n  [\"338\u00d76=2028\", \"588\u00d76=3528\"],\n  [\"289\u00d73=867\", \"798\u00d75=3990\"],\n  [\"201\u00d79=1809\", \"511\u00d76=3066\"],\n  [\"713\u00d74=2852\", \"984\u00d73=2952\"],\n  [\"239\u00d77=1673\", \"345\u00d74=1380\"],\n  [\"561\u00d75=2805\", \"147\u00d79=1323\"],\n  [\"721\u00d76=4326\", \"367\u00d72=734\"],\n  [\"794\u00d79=7146\", \"464\u00d76=2784\"],\n  [\"528\u00d72=1056\", \"255\u00d76=1530\"],\n  [\"212\u00d74=848\", \"264\u00d75=1320\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw `Could not find text: ${oldText}`;\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit x one-digit multiplication expression with its\n# corresponding updated expression, one Find/Replace pair per cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"581\u00d73=1743\"; New = \"530\u00d73=1590\" }\n    @{ Old = \"551\u00d73=1653\"; New = \"639\u00d79=5751\" }\n    @{ Old = \"411\u00d72=822\"; New = \"514\u00d74=2056\" }\n    @{ Old = \"789\u00d77=5523\"; New = \"124\u00d73=372\" }\n    @{ Old = \"960\u00d75=4800\"; New = \"469\u00d79=4221\" }\n    @{ Old = \"872\u00d79=7848\"; New = \"806\u00d79=7254\" }\n    @{ Old = \"263\u00d77=1841\"; New = \"978\u00d77=6846\" }\n    @{ Old = \"716\u00d77=5012\"; New = \"278\u00d76=1668\" }\n    @{ Old = \"369\u00d78=2952\"; New = \"343\u00d74=1372\" }\n    @{ Old = \"980\u00d78=7840\"; New = \"336\u00d72=672\" }\n    @{ Old = \"389\u00d78=3112\"; New = \"245\u00d74=980\" }\n    @{ Old = \"346\u00d74=1384\"; New = \"392\u00d72=784\" }\n    @{ Old = \"538\u00d74=2152\"; New = \"825\u00d72=1650\" }\n    @{ Old = \"179\u00d75=895\"; New = \"974\u00d72=1948\" }\n    @{ Old = \"554\u00d73=1662\"; New = \"604\u00d72=1208\" }\n    @{ Old = \"338\u00d76=2028\"; New = \"588\u00d76=3528\" }\n    @{ Old = \"289\u00d73=867\"; New = \"798\u00d75=3990\" }\n    @{ Old = \"201\u00d79=1809\"; New = \"511\u00d76=3066\" }\n    @{ Old = \"713\u00d74=2852\"; New = \"984\u00d73=2952\" }\n    @{ Old = \"239\u00d77=1673\"; New = \"345\u00d74=1380\" }\n    @{ Old = \"561\u00d75=2805\"; New = \"147\u00d79=1323\" }\n    @{ Old = \"721\u00d76=4326\"; New = \"367\u00d72=734\" }\n    @{ Old = \"794\u00d79=7146\"; New = \"464\u00d76=2784\" }\n    @{ Old = \"528\u00d72=1056\"; New = \"255\u00d76=1530\" }\n    @{ Old = \"212\u00d74=848\"; New = \"264\u00d75=1320\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
